$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45085
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 19000
$ws.Range("P2").Value = 19000
$ws.Range("S2").Value = 1056

# Row 3
$ws.Range("D3").Value = 45069
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("S3").Value = 833

# Row 4
$ws.Range("D4").Value = 45096
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("S4").Value = 1111

# Row 5
$ws.Range("D5").Value = 45076
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("S5").Value = 833

# Row 6
$ws.Range("D6").Value = 45084
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 19000
$ws.Range("P6").Value = 18500
$ws.Range("S6").Value = 1028

# Row 7
$ws.Range("D7").Value = 45112
$ws.Range("M7").Value = 20
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("S7").Value = 1111

# Row 8
$ws.Range("D8").Value = 45111
$ws.Range("M8").Value = 20

# Row 10
$ws.Range("D10").Value = 45072
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("S10").Value = 833

# Row 11
$ws.Range("D11").Value = 45061
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("S11").Value = 833

# Row 12
$ws.Range("D12").Value = 45092
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 19000
$ws.Range("P12").Value = 18667
$ws.Range("S12").Value = 1037

# Row 13
$ws.Range("D13").Value = 45083
$ws.Range("M13").Value = 50
